$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set Quantity values for the camber-adjustment shim rows
$ws.Range("F5").Value = 6
$ws.Range("F12").Value = 6

# Update the remembered selection/active cell on the sheet
$ws.Range("E17").Select()
